# Update "想去人数" (interested-attendee count) figures on both the
# "展览" and "全部类型" sheets to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

# Sheet "展览": rows 4, 5, 8 in column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 175
$wsExhibit.Range("F5").Value = 3294
$wsExhibit.Range("F8").Value = 417

# Sheet "全部类型": rows 4, 5, 10 in column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 175
$wsAll.Range("F5").Value = 3294
$wsAll.Range("F10").Value = 417
